$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the C6 value (was a duplicate/typo "9991231231239" -> "9999999999999")
$ws.Range("C6").Value = "9999999999999"

# Equivalence Class Partitioning (ECP) block
$ws.Range("A10").Value = "Equivalence Class Partitioning (ECP):"
$ws.Range("A11").Value = "Gültige Eingaben:"
$ws.Range("A12").Value = "CNP-Nummern mit gültigen Strukturen und Werten."
$ws.Range("A13").Value = "CNP-Nummern mit unterschiedlichen Gültigkeitsprüfungen."

$ws.Range("A15").Value = "Ungültige Eingaben:"

# Boundary Value Analysis (BVA) block
$ws.Range("A19").Value = "Boundary Value Analysis (BVA):"
$ws.Range("A20").Value = "Gültige Grenzwerte:"
$ws.Range("A21").Value = "Die kleinste gültige CNP-Nummer."
$ws.Range("A22").Value = "Die größte gültige CNP-Nummer."

$ws.Range("A24").Value = "Ungültige Grenzwerte:"
$ws.Range("A26").Value = "bzw. größten gültigen CNP liegen."

# Back-fill the two invalid-ECP lines, then the invalid-BVA wrap line
$ws.Range("A16").Value = "CNP-Nummern mit ungültigen Strukturen (falsche Länge, ungültige Zeichen usw.)."
$ws.Range("A17").Value = "CNP-Nummern mit ungültigen Werten (falsche Prüfsumme, ungültige Geburtsdaten usw.)."
$ws.Range("A25").Value = "CNP-Nummern, die eine Stelle unter oder über dem kleinsten"

# Column A needs to widen to fit the longest new text
$ws.Columns("A").AutoFit() | Out-Null

# Move selection/active cell to the row after the new content, like Excel leaves after typing
$ws.Range("A27").Select() | Out-Null
